# Auto-generated Excel COM-interop script to apply the Siren_Profits market-data refresh.
# Updates cached price/profit figures (columns H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 252.57143
$ws.Range("I33").Value = 256.6154
$ws.Range("K33").Value = 256.6154
$ws.Range("M33").Value = -27.61540000000002
# Row 70
$ws.Range("H70").Value = 4590.636
$ws.Range("I70").Value = 6874.5
$ws.Range("J70").Value = 3285.5715
$ws.Range("K70").Value = 20623.5
$ws.Range("L70").Value = 9856.7145
$ws.Range("M70").Value = -20353.5
$ws.Range("N70").Value = -10396.7145
# Row 73
$ws.Range("H73").Value = 4590.636
$ws.Range("I73").Value = 6874.5
$ws.Range("J73").Value = 3285.5715
$ws.Range("K73").Value = 20623.5
$ws.Range("L73").Value = 9856.7145
$ws.Range("M73").Value = -19687.5
$ws.Range("N73").Value = -11728.7145
# Row 74
$ws.Range("H74").Value = 4489.778
$ws.Range("I74").Value = 3865.2
$ws.Range("K74").Value = 3865.2
$ws.Range("M74").Value = -2929.2
# Row 77
$ws.Range("H77").Value = 4489.778
$ws.Range("I77").Value = 3865.2
$ws.Range("K77").Value = 19326
$ws.Range("M77").Value = -14646
# Row 100
$ws.Range("H100").Value = 7036598
$ws.Range("I100").Value = 8743.111000000001
$ws.Range("K100").Value = 8743.111000000001
$ws.Range("M100").Value = -8202.111000000001
# Row 138
$ws.Range("H138").Value = 4393.2705
$ws.Range("J138").Value = 5179.826
$ws.Range("L138").Value = 15539.478
$ws.Range("N138").Value = -25819.478

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 9012.92
$ws.Range("I45").Value = 9372.263000000001
$ws.Range("K45").Value = 9372.263000000001
$ws.Range("M45").Value = -8995.263000000001
# Row 61
$ws.Range("H61").Value = 6379.5557
$ws.Range("I61").Value = 6166.5884
$ws.Range("K61").Value = 6166.5884
$ws.Range("M61").Value = -5954.5884
# Row 136
$ws.Range("H136").Value = 6379.5557
$ws.Range("I136").Value = 6166.5884
$ws.Range("K136").Value = 18499.7652
$ws.Range("M136").Value = -15949.7652

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 97
$ws.Range("H97").Value = 19997.666
$ws.Range("I97").Value = 14996.5
$ws.Range("J97").Value = 30000
$ws.Range("K97").Value = 14996.5
$ws.Range("L97").Value = 30000
$ws.Range("M97").Value = -14005.5
$ws.Range("N97").Value = -31982
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
# Row 141
$ws.Range("H141").Value = 68587.60000000001
$ws.Range("J141").Value = 68587.60000000001
$ws.Range("L141").Value = 68587.60000000001
$ws.Range("N141").Value = -78947.60000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 7368.4814
$ws.Range("I58").Value = 15074.111
$ws.Range("K58").Value = 15074.111
$ws.Range("M58").Value = -14871.111
# Row 86
$ws.Range("H86").Value = 11296.833
$ws.Range("I86").Value = 9882
$ws.Range("J86").Value = 16248.75
$ws.Range("K86").Value = 9882
$ws.Range("L86").Value = 16248.75
$ws.Range("M86").Value = -8759
$ws.Range("N86").Value = -18494.75
# Row 89
$ws.Range("H89").Value = 11296.833
$ws.Range("I89").Value = 9882
$ws.Range("J89").Value = 16248.75
$ws.Range("K89").Value = 49410
$ws.Range("L89").Value = 81243.75
$ws.Range("M89").Value = -43794
$ws.Range("N89").Value = -92475.75
# Row 134
$ws.Range("H134").Value = 2238394.8
$ws.Range("I134").Value = 2848143
$ws.Range("K134").Value = 8544429
$ws.Range("M134").Value = -8541894
# Row 136
$ws.Range("H136").Value = 7368.4814
$ws.Range("I136").Value = 15074.111
$ws.Range("K136").Value = 45222.333
$ws.Range("M136").Value = -42672.333

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 716.4474
$ws.Range("I86").Value = 950
$ws.Range("J86").Value = 681.0606
$ws.Range("K86").Value = 2850
$ws.Range("L86").Value = 2043.1818
$ws.Range("M86").Value = -1664
$ws.Range("N86").Value = -4415.1818
# Row 89
$ws.Range("H89").Value = 716.4474
$ws.Range("I89").Value = 950
$ws.Range("J89").Value = 681.0606
$ws.Range("K89").Value = 8550
$ws.Range("L89").Value = 6129.5454
$ws.Range("M89").Value = -2622
$ws.Range("N89").Value = -17985.5454
# Row 97
$ws.Range("H97").Value = 36577.35
$ws.Range("I97").Value = 86582.14
$ws.Range("J97").Value = 1574
$ws.Range("K97").Value = 259746.42
$ws.Range("L97").Value = 4722
$ws.Range("M97").Value = -259250.42
$ws.Range("N97").Value = -5714
# Row 98
$ws.Range("H98").Value = 762.6667
$ws.Range("I98").Value = 947.3333
$ws.Range("J98").Value = 393.33334
$ws.Range("K98").Value = 2841.9999
$ws.Range("L98").Value = 1180.00002
$ws.Range("M98").Value = -1343.9999
$ws.Range("N98").Value = -4176.000019999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 250.375
$ws.Range("J2").Value = 243.9
$ws.Range("L2").Value = 243.9
$ws.Range("N2").Value = -469.9
# Row 70
$ws.Range("H70").Value = 9950.700000000001
$ws.Range("I70").Value = 9001
$ws.Range("K70").Value = 9001
$ws.Range("M70").Value = -8731
# Row 73
$ws.Range("H73").Value = 9950.700000000001
$ws.Range("I73").Value = 9001
$ws.Range("K73").Value = 9001
$ws.Range("M73").Value = -8065
# Row 101
$ws.Range("H101").Value = 41999
$ws.Range("I101").Value = 30000
$ws.Range("J101").Value = 44998.75
$ws.Range("K101").Value = 30000
$ws.Range("L101").Value = 44998.75
$ws.Range("M101").Value = -26755
$ws.Range("N101").Value = -51488.75
# Row 126
$ws.Range("H126").Value = 27728.545
$ws.Range("I126").Value = 44002.5
$ws.Range("K126").Value = 132007.5
$ws.Range("M126").Value = -129537.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 621.3022999999999
$ws.Range("I55").Value = 620.7406999999999
$ws.Range("J55").Value = 622.25
$ws.Range("K55").Value = 620.7406999999999
$ws.Range("L55").Value = 622.25
$ws.Range("M55").Value = -447.7406999999999
$ws.Range("N55").Value = -968.25
# Row 61
$ws.Range("H61").Value = 2987.5
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 3650
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 3650
$ws.Range("M61").Value = -798
$ws.Range("N61").Value = -4054
# Row 113
$ws.Range("H113").Value = 2987.5
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 3650
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 3650
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -7990
# Row 136
$ws.Range("H136").Value = 8550.261
$ws.Range("I136").Value = 3337.25
$ws.Range("J136").Value = 11330.533
$ws.Range("K136").Value = 10011.75
$ws.Range("L136").Value = 33991.599
$ws.Range("M136").Value = -7461.75
$ws.Range("N136").Value = -39091.599

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2800.875
$ws.Range("J96").Value = 3620.8
$ws.Range("L96").Value = 3620.8
$ws.Range("N96").Value = -6366.8
# Row 132
$ws.Range("H132").Value = 31393.908
$ws.Range("I132").Value = 34233.8
$ws.Range("K132").Value = 102701.4
$ws.Range("M132").Value = -100171.4
# Row 136
$ws.Range("H136").Value = 1883.4062
$ws.Range("I136").Value = 1533.6923
$ws.Range("J136").Value = 3398.8333
$ws.Range("K136").Value = 4601.0769
$ws.Range("L136").Value = 10196.4999
$ws.Range("M136").Value = -2051.0769
$ws.Range("N136").Value = -15296.4999
